$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-20 11:01:44", 0.0004),
    @("2023-12-20 11:02:15", 0.0018),
    @("2023-12-20 11:03:13", 0.003800000000000001),
    @("2023-12-20 11:03:18", 0.0004),
    @("2023-12-20 11:03:22", 0.0004),
    @("2023-12-20 11:03:57", 0.001)
)

$startRow = 511
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
